$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bug")

# ---------------------------------------------------------------------------
# New / updated bug rows 38-43
# Columns: B=No. C=id D=found platform E=found date F=description G=status H=root cause
# ---------------------------------------------------------------------------

# Rows 38-41 and 43 carry the "yellow highlight" fill (same as row 37); row 42
# stays on the plain (no fill) style.
$ws.Range("B38:H38").Interior.Color = 65535
$ws.Range("B39:H39").Interior.Color = 65535
$ws.Range("B40:H40").Interior.Color = 65535
$ws.Range("B41:H41").Interior.Color = 65535
$ws.Range("B43:H43").Interior.Color = 65535

$ws.Cells.Item(38, 2).Value = 37
$ws.Cells.Item(38, 3).Value = "PPBOX-3717 "
$ws.Cells.Item(38, 4).Value = "OTT"
$ws.Cells.Item(38, 5).Value = 20150320
$ws.Cells.Item(38, 6).Value = "DLNA推送在线视频，视频A播放完后播放视频B时，提示视频解析出错，实际可正常播放"
$ws.Cells.Item(38, 7).Value = "tracking"

$ws.Cells.Item(39, 2).Value = 38
$ws.Cells.Item(39, 3).Value = "PPBOX-3706"
$ws.Cells.Item(39, 4).Value = "OTT"
$ws.Cells.Item(39, 5).Value = 20150320
$ws.Cells.Item(39, 6).Value = "本地播放.ogg音乐，不论播放控制选择什么，播放模式都是单个循环"
$ws.Cells.Item(39, 7).Value = "tracking"
$ws.Cells.Item(39, 8).Value = "系统播放器 没有 onComplete 回调？"

$ws.Cells.Item(40, 2).Value = 39
$ws.Cells.Item(40, 3).Value = "PPBOX-3693"
$ws.Cells.Item(40, 4).Value = "OTT"
$ws.Cells.Item(40, 5).Value = 20150320
$ws.Cells.Item(40, 6).Value = "外挂字幕的片源，播放时不会自动挂载字幕，手动加载也加载不上"
$ws.Cells.Item(40, 7).Value = "tracking"
$ws.Cells.Item(40, 8).Value = "老版sdk支持的srt字幕，现在不支持了"

$ws.Cells.Item(41, 2).Value = 40
$ws.Cells.Item(41, 3).Value = "PPBOX-3417"
$ws.Cells.Item(41, 4).Value = "OTT"
$ws.Cells.Item(41, 5).Value = 20150320
$ws.Cells.Item(41, 6).Value = "wmv片源新sdk版本播放声音卡顿，旧sdk正常"
$ws.Cells.Item(41, 7).Value = "tracking"

$ws.Cells.Item(42, 2).Value = 41
$ws.Cells.Item(42, 3).Value = "PPBOX-3695"
$ws.Cells.Item(42, 4).Value = "OTT"
$ws.Cells.Item(42, 5).Value = 20150319
$ws.Cells.Item(42, 6).Value = "播放多音轨片源 选择音轨2时，我的设备挂掉"
$ws.Cells.Item(42, 7).Value = "fixed"
$ws.Cells.Item(42, 8).Value = "多音轨 各个音轨的channel_layout属性不同，需要重启audioplayer"

$ws.Cells.Item(43, 2).Value = 42
$ws.Cells.Item(43, 3).Value = "N/A"
$ws.Cells.Item(43, 4).Value = "IOS"
$ws.Cells.Item(43, 5).Value = 20150318
$ws.Cells.Item(43, 6).Value = "IOS6系统 自有播放器 播放hls crash"
$ws.Cells.Item(43, 7).Value = "tracking"
$ws.Cells.Item(43, 8).Value = "ld error"

# ---------------------------------------------------------------------------
# Trailing placeholder rows 44-52 - just a continuing No. sequence in column B
# ---------------------------------------------------------------------------
for ($i = 44; $i -le 52; $i++) {
    $ws.Cells.Item($i, 2).Value = $i - 1
}

# ---------------------------------------------------------------------------
# View state: scrolled so row 22 is at the top, selection moved to D32
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("D32").Select()
